# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml) - rows keyed by row number, new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2
$ws1.Range("F5").Value = 16
$ws1.Range("F7").Value = 1643
$ws1.Range("F8").Value = 6
$ws1.Range("F9").Value = 10
$ws1.Range("F11").Value = 1493
$ws1.Range("F14").Value = 375
$ws1.Range("F15").Value = 251
$ws1.Range("F16").Value = 188
$ws1.Range("F19").Value = 21
$ws1.Range("F21").Value = 266
$ws1.Range("F22").Value = 146
$ws1.Range("F23").Value = 213
$ws1.Range("F24").Value = 203

# Sheet "全部类型" (sheet4.xml) - same events, rows shifted by +1 after row 7
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2
$ws4.Range("F5").Value = 16
$ws4.Range("F7").Value = 1643
$ws4.Range("F9").Value = 6
$ws4.Range("F10").Value = 10
$ws4.Range("F12").Value = 1493
$ws4.Range("F15").Value = 375
$ws4.Range("F16").Value = 251
$ws4.Range("F17").Value = 188
$ws4.Range("F20").Value = 21
$ws4.Range("F22").Value = 266
$ws4.Range("F23").Value = 146
$ws4.Range("F24").Value = 213
$ws4.Range("F25").Value = 203
